$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.515.49'
$ws.Range('E2').Value = '  +6.15%  '
$ws.Range('D3').Value = '3.693.23'
$ws.Range('E3').Value = '  +6.02%  '
$ws.Range('E4').Value = '  +0.67%  '
$ws.Range('D5').Value = "'425.71"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.55%  '
$ws.Range('D6').Value = "'130.34"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.43%  '
$ws.Range('D7').Value = '3.675.93'
$ws.Range('E7').Value = '  +5.76%  '
$ws.Range('D8').Value = "'0.639"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.21%  '
$ws.Range('E9').Value = '  +0.08%  '
$ws.Range('D10').Value = "'0.764"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.45%  '
$ws.Range('D11').Value = "'0.179"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +12.75%  '
$ws.Range('D12').Value = "'0.0000361"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +58.21%  '
$ws.Range('D13').Value = "'42.33"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.27%  '
$ws.Range('D14').Value = "'9.94"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.74%  '
$ws.Range('D15').Value = '4.286.03'
$ws.Range('E15').Value = '  +6.14%  '
$ws.Range('E16').Value = '  +0.09%  '
$ws.Range('D17').Value = "'20.49"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.90%  '
$ws.Range('D18').Value = '3.695.55'
$ws.Range('E18').Value = '  +5.97%  '
$ws.Range('D19').Value = "'1.12"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.80%  '
$ws.Range('D20').Value = "'12.73"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.29%  '
$ws.Range('D21').Value = '67.549.93'
$ws.Range('E21').Value = '  +6.23%  '
$ws.Range('D22').Value = "'446.66"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.56%  '
$ws.Range('D23').Value = "'15.37"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +16.20%  '
$ws.Range('D24').Value = "'88.94"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.56%  '
$ws.Range('D25').Value = "'3.11"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -5.36%  '
$ws.Range('D26').Value = "'37.32"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +11.29%  '
$ws.Range('D27').Value = "'10.42"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.34%  '
$ws.Range('D28').Value = "'3.30"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.25%  '
$ws.Range('E29').Value = '  +4.33%  '
$ws.Range('D30').Value = "'2.78"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +3.50%  '
$ws.Range('D31').Value = "'12.43"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.32%  '
$ws.Range('E32').Value = '  +7.24%  '
$ws.Range('D33').Value = "'7.17"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -4.57%  '
$ws.Range('D34').Value = "'0.160"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.12%  '
$ws.Range('D35').Value = "'40.46"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.93%  '
$ws.Range('D36').Value = "'0.998"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.10%  '
$ws.Range('D37').Value = "'56.34"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.18%  '
$ws.Range('D38').Value = "'0.0489"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.59%  '
$ws.Range('D39').Value = '0.0₃0715'
$ws.Range('E39').Value = '  +10.70%  '
$ws.Range('B40').Value = 'Stellar'
$ws.Range('C40').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D40').Value = "'0.146"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +6.14%  '
$ws.Range('B41').Value = 'ThetaToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D41').Value = "'2.93"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +25.01%  '
$ws.Range('E42').Value = '  -0.02%  '
$ws.Range('D43').Value = "'3.41"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.98%  '
$ws.Range('D44').Value = "'146.95"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.51%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').Value = "'26.91"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +23.02%  '
$ws.Range('B46').Value = 'Stacks'
$ws.Range('C46').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D46').Value = "'2.92"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -5.78%  '
$ws.Range('B47').Value = 'WEMIXToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D47').Value = "'2.67"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -5.36%  '
$ws.Range('D48').Value = "'2.07"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.30%  '
$ws.Range('D49').Value = "'4.31"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -4.61%  '
$ws.Range('D50').Value = "'0.303"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.22%  '
$ws.Range('D51').Value = "'0.158"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +13.12%  '
